# Refresh crypto price/volume snapshot (GitHub Actions run, 2023-10-28 23:27 UTC).
# Row 13/14 additionally swap Chainlink <-> WrappedEther (ranking reshuffle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.154.69'
$ws.Range("D3").Value = '1.777.90'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  +0.12%  '
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '225.81'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +0.15%  '
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '31.63'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("E9").Value = '  +0.87%  '
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0690'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +1.87%  '
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("D12").Value = '2.034.29'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.776.26'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.93'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("D15").Value = '34.111.35'
$ws.Range("E15").Value = '  +0.75%  '
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.621'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("E17").Value = '  +1.09%  '
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '67.81'
$cell.Style = $origStyle
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("D19").Value = '0.0₃0798'
$ws.Range("E19").Value = '  +3.62%  '
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '245.40'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +2.64%  '
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.96'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +3.88%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  +1.71%  '
$ws.Range("E24").Value = '  -0.87%  '
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '162.52'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("E26").Value = '  +2.37%  '
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '16.28'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("E28").Value = '  +1.72%  '
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("E32").Value = '  +3.90%  '
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.70'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +5.20%  '
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = '1.437.98'
$ws.Range("E35").Value = '  +3.42%  '
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.661'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +3.84%  '
$ws.Range("E37").Value = '  +6.52%  '
$ws.Range("E38").Value = '  +2.29%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("E41").Value = '  -0.15%  '
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.920'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +0.66%  '
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.66'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +0.60%  '
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.44'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("E45").Value = '  +0.20%  '
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.07'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +3.62%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").Value = '0.0₆0135'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '1.936.28'
$ws.Range("E49").Value = '  -0.10%  '
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '104.13'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("E51").Value = '  +0.17%  '
